{"js": "// Replace the 25 multiplication-equation answers in the single body table\n// with their updated values, in document (row-major) order. Blank cells\n// (the spacer rows) are left untouched. This mirrors the diff exactly:\n// same table shape (20 rows x 5 cols), only the <w:t> text of the\n// non-empty cells changes.\n\nconst newValues = [\n  \"116\u00d77=812\", \"932\u00d74=3728\", \"614\u00d75=3070\", \"581\u00d76=3486\", \"231\u00d76=1386\",\n  \"629\u00d79=5661\", \"887\u00d75=4435\", \"671\u00d76=4026\", \"374\u00d77=2618\", \"261\u00d73=783\",\n  \"396\u00d77=2772\", \"809\u00d75=4045\", \"429\u00d78=3432\", \"663\u00d79=5967\", \"808\u00d77=5656\",\n  \"295\u00d78=2360\", \"631\u00d76=3786\", \"900\u00d74=3600\", \"101\u00d76=606\", \"861\u00d78=6888\",\n  \"502\u00d76=3012\", \"160\u00d75=800\", \"330\u00d75=1650\", \"735\u00d76=4410\", \"866\u00d76=5196\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values;\nlet i = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    if (grid[r][c] !== \"\") {\n      grid[r][c] = newValues[i];\n      i++;\n    }\n  }\n}\n\ntable.values = grid;\nawait context.sync();\n", "ps1": "# Replace the 25 multiplication-equation answers in the single body table\n# with their updated values, in document (row-major) order. Blank spacer\n# cells are left untouched. Mirrors the diff exactly: same table shape\n# (20 rows x 5 cols), only the text of the non-empty cells changes.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n  \"116\u00d77=812\",\n  \"932\u00d74=3728\",\n  \"614\u00d75=3070\",\n  \"581\u00d76=3486\",\n  \"231\u00d76=1386\",\n  \"629\u00d79=5661\",\n  \"887\u00d75=4435\",\n  \"671\u00d76=4026\",\n  \"374\u00d77=2618\",\n  \"261\u00d73=783\",\n  \"396\u00d77=2772\",\n  \"809\u00d75=4045\",\n  \"429\u00d78=3432\",\n  \"663\u00d79=5967\",\n  \"808\u00d77=5656\",\n  \"295\u00d78=2360\",\n  \"631\u00d76=3786\",\n  \"900\u00d74=3600\",\n  \"101\u00d76=606\",\n  \"861\u00d78=6888\",\n  \"502\u00d76=3012\",\n  \"160\u00d75=800\",\n  \"330\u00d75=1650\",\n  \"735\u00d76=4410\",\n  \"866\u00d76=5196\"\n)\n\n$idx = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n  for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n    $cell = $tbl.Cell($r, $c)\n    $txt = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -ne \"\") {\n      $cell.Range.Text = $newValues[$idx]\n      $idx++\n    }\n  }\n}\n"}
